$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.279.47'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +4.74%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.453.86'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.30%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.14%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.31'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.65%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '183.62'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +6.15%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.53%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.447.27'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +3.25%  '

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.01%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.177'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +9.15%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.646'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +3.46%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '55.67'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +4.22%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000281'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.20%  '

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +4.23%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.006.28'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.94%  '

$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.53'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.48%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.446.04'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.13%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '67.369.02'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +4.82%  '

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.41%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.03'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.36%  '

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.09%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '485.49'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +6.23%  '

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.99%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '15.12'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +12.22%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.19'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.71%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '89.72'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +5.00%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.97'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.96%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.92'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.33%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.92'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.21%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.59'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +4.29%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.96'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +5.00%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '600.74'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +5.32%  '

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.28%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.00'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.88%  '

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.94%  '

$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.148'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +6.79%  '

$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.04%  '

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.44%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0₃0782'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +6.58%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.389'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +6.28%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.50'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +4.17%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.152.68'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.04%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.93'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +5.19%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.59'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +6.70%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0426'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.53%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.82'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +23.54%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.29'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +4.89%  '

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.53%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.75'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +8.08%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.01%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '141.35'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.42%  '
